$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.914.60"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.228.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.70"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.98"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.77"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.97"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.584.36"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.04"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.239.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.728"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.897.57"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.79"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.29%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.75"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.22"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.88"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0715"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.22"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.45%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0997"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.70"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.050.40"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.03"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +13.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.454.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.50"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.98"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.09%  "
